$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("devices")

# The capabilities row was missing the automation name value; add it.
$ws.Range("K2").Value = "PerfectoMobile"

# Update the active selection to the cell we just edited.
$ws.Range("K2").Select()
